$wb = $excel.ActiveWorkbook

# --- About sheet: insert a new "note" block before the existing time-series note ---
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("55:57").Insert()
$wsAbout.Range("A55").Value = "The units for rail in the BTS data set are unclear whether they report locomotive-miles or train-car-miles. Based on the"
$wsAbout.Range("A56").Value = "scale, we assume the units are train-car-miles."

# --- BTS NTS Modal Profile Data sheet: insert new assumption row for train cars per locomotive ---
$wsBts = $wb.Worksheets.Item("BTS NTS Modal Profile Data")
$wsBts.Range("22:22").Insert()
$wsBts.Range("A22").Value = "Assumption - train cars per locomotive"
$wsBts.Range("B22").Value = 10
$wsBts.Range("A22:B22").Font.Bold = $true
$wsBts.Range("A22:B22").Interior.Pattern = -4142

# Update the weighted-average passenger rail loading formula (now row 37 after the insert)
# to scale per-locomotive loading up to a per-train basis using the new assumption cell.
$wsBts.Range("B37").Formula = "=(B26*B25+B34*B28+B35*B29+B36*B30)/SUM(B25,B28:B30)*B22"
